# Generate Report for Handoff
# Adds two newly-handed-off files (4a674e25-... and 878fd696-...) to the
# localization status report: one new row each on the Overview sheet and
# on the per-language (zh-cn / de-de) detail sheets, pushing the existing
# ".localization-config" row down by two rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Shared URL fragments used to rebuild hyperlinks.
# ---------------------------------------------------------------------
$mdBase     = "https://github.com/OpenLocalizationTest/oltest/blob/041cd0bce2ae55136360588cd7a6585f43442cf6/e2e/"
$configUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/041cd0bce2ae55136360588cd7a6585f43442cf6/.localization-config"
$zhXlfBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7c63bc377c71a8525dddb5735dfbd23c36129465/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/"
$deXlfBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dfd61c89e4eb849ed0b082c6f78e03cb463c6238/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/"

$file1 = "4a674e25-3ba2-4e52-833b-68918e322936.md"
$file2 = "878fd696-c78d-4955-9f93-3473f0657199.md"

$file1ZhXlf = "4a674e25-3ba2-4e52-833b-68918e322936.f470f54e3993e463d5a8ee8c1e80c4044f0f6bd7.zh-cn.xlf"
$file2ZhXlf = "878fd696-c78d-4955-9f93-3473f0657199.b770f6653a0056dfab75e6a68c30d6d72a99a4fc.zh-cn.xlf"
$file1DeXlf = "4a674e25-3ba2-4e52-833b-68918e322936.f470f54e3993e463d5a8ee8c1e80c4044f0f6bd7.de-de.xlf"
$file2DeXlf = "878fd696-c78d-4955-9f93-3473f0657199.b770f6653a0056dfab75e6a68c30d6d72a99a4fc.de-de.xlf"

$zhDate = "2016-02-24 07:33:09"
$deDate = "2016-02-24 07:33:20"
$epoch  = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# 1) Overview sheet — columns A:C, rows 1..4 -> 1..6
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Rows.Item(4).Insert()
$ov.Rows.Item(4).Insert()

$ov.Range("A4").Value = $file1
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"

$ov.Range("A5").Value = $file2
$ov.Range("B5").Value = "Ready for handoff"
$ov.Range("C5").Value = "Ready for handoff"

# Row 6 already carries the shifted ".localization-config" / "Not to be
# localized" values forward from the old row 4 thanks to Rows.Insert().

# Hyperlink `ref`s are not re-anchored by Rows.Insert(), so rebuild the
# whole collection against the final row layout.
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), $mdBase + "6d2e9fb0-015e-48a4-991f-2a6b03475b8e.md", "", "", "6d2e9fb0-015e-48a4-991f-2a6b03475b8e.md")
$ov.Hyperlinks.Add($ov.Range("A3"), $mdBase + "fe0528f8-4d62-44e2-a2a7-ea8e1baf428f.md", "", "", "fe0528f8-4d62-44e2-a2a7-ea8e1baf428f.md")
$ov.Hyperlinks.Add($ov.Range("A4"), $mdBase + $file1, "", "", $file1)
$ov.Hyperlinks.Add($ov.Range("A5"), $mdBase + $file2, "", "", $file2)
$ov.Hyperlinks.Add($ov.Range("A6"), $configUrl, "", "", ".localization-config")

# ---------------------------------------------------------------------
# 2) zh-cn sheet — columns A:I, rows 1..4 -> 1..6
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Rows.Item(4).Insert()
$zh.Rows.Item(4).Insert()

$zh.Range("A4").Value = $file1
$zh.Range("B4").Value = "Ready for handoff"
$zh.Range("C4").Value = $file1ZhXlf
$zh.Range("D4").Value = $zhDate
$zh.Range("G4").Value = $epoch
$zh.Range("H4").Value = "Include"

$zh.Range("A5").Value = $file2
$zh.Range("B5").Value = "Ready for handoff"
$zh.Range("C5").Value = $file2ZhXlf
$zh.Range("D5").Value = $zhDate
$zh.Range("G5").Value = $epoch
$zh.Range("H5").Value = "Include"

# Row 6 (shifted ".localization-config" row) keeps its old D/G/H values;
# only make sure the formerly-empty C column stays empty (it already is).

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), $mdBase + "6d2e9fb0-015e-48a4-991f-2a6b03475b8e.md", "", "", "6d2e9fb0-015e-48a4-991f-2a6b03475b8e.md")
$zh.Hyperlinks.Add($zh.Range("C2"), $zhXlfBase + "6d2e9fb0-015e-48a4-991f-2a6b03475b8e.9b44a5e302e2c98dd79ec253ed6cb9040a68128e.zh-cn.xlf", "", "", "6d2e9fb0-015e-48a4-991f-2a6b03475b8e.9b44a5e302e2c98dd79ec253ed6cb9040a68128e.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A3"), $mdBase + "fe0528f8-4d62-44e2-a2a7-ea8e1baf428f.md", "", "", "fe0528f8-4d62-44e2-a2a7-ea8e1baf428f.md")
$zh.Hyperlinks.Add($zh.Range("C3"), $zhXlfBase + "fe0528f8-4d62-44e2-a2a7-ea8e1baf428f.70460a93df7465fe6297039e43cad0efe0da6720.zh-cn.xlf", "", "", "fe0528f8-4d62-44e2-a2a7-ea8e1baf428f.70460a93df7465fe6297039e43cad0efe0da6720.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A4"), $mdBase + $file1, "", "", $file1)
$zh.Hyperlinks.Add($zh.Range("C4"), $zhXlfBase + $file1ZhXlf, "", "", $file1ZhXlf)
$zh.Hyperlinks.Add($zh.Range("A5"), $mdBase + $file2, "", "", $file2)
$zh.Hyperlinks.Add($zh.Range("C5"), $zhXlfBase + $file2ZhXlf, "", "", $file2ZhXlf)
$zh.Hyperlinks.Add($zh.Range("A6"), $configUrl, "", "", ".localization-config")

# ---------------------------------------------------------------------
# 3) de-de sheet — columns A:I, rows 1..4 -> 1..6
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Rows.Item(4).Insert()
$de.Rows.Item(4).Insert()

$de.Range("A4").Value = $file1
$de.Range("B4").Value = "Ready for handoff"
$de.Range("C4").Value = $file1DeXlf
$de.Range("D4").Value = $deDate
$de.Range("G4").Value = $epoch
$de.Range("H4").Value = "Include"

$de.Range("A5").Value = $file2
$de.Range("B5").Value = "Ready for handoff"
$de.Range("C5").Value = $file2DeXlf
$de.Range("D5").Value = $deDate
$de.Range("G5").Value = $epoch
$de.Range("H5").Value = "Include"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), $mdBase + "6d2e9fb0-015e-48a4-991f-2a6b03475b8e.md", "", "", "6d2e9fb0-015e-48a4-991f-2a6b03475b8e.md")
$de.Hyperlinks.Add($de.Range("C2"), $deXlfBase + "6d2e9fb0-015e-48a4-991f-2a6b03475b8e.9b44a5e302e2c98dd79ec253ed6cb9040a68128e.de-de.xlf", "", "", "6d2e9fb0-015e-48a4-991f-2a6b03475b8e.9b44a5e302e2c98dd79ec253ed6cb9040a68128e.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A3"), $mdBase + "fe0528f8-4d62-44e2-a2a7-ea8e1baf428f.md", "", "", "fe0528f8-4d62-44e2-a2a7-ea8e1baf428f.md")
$de.Hyperlinks.Add($de.Range("C3"), $deXlfBase + "fe0528f8-4d62-44e2-a2a7-ea8e1baf428f.70460a93df7465fe6297039e43cad0efe0da6720.de-de.xlf", "", "", "fe0528f8-4d62-44e2-a2a7-ea8e1baf428f.70460a93df7465fe6297039e43cad0efe0da6720.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A4"), $mdBase + $file1, "", "", $file1)
$de.Hyperlinks.Add($de.Range("C4"), $deXlfBase + $file1DeXlf, "", "", $file1DeXlf)
$de.Hyperlinks.Add($de.Range("A5"), $mdBase + $file2, "", "", $file2)
$de.Hyperlinks.Add($de.Range("C5"), $deXlfBase + $file2DeXlf, "", "", $file2DeXlf)
$de.Hyperlinks.Add($de.Range("A6"), $configUrl, "", "", ".localization-config")
